# Underwater_Robot_NFC balance sheet update:
#  - "No.12" row (B12) balance was corrected from 490 back to 500
#  - "No.14" row (B14) balance was corrected from 495 back to 500
#  - the active selection in the sheet view was moved from D15 to B7:B17
#    (reviewing/confirming the balance column so Excel updates show live)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 500
$ws.Range("B14").Value = 500

$ws.Range("B7:B17").Select()
